$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.365.17"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.466.15"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.604"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.72%  "
$ws.Range("D8").Value = "3.458.68"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.570"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000274"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "4.038.22"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "613.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.58%  "
$ws.Range("D18").Value = "3.478.11"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "68.488.42"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.868"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "570.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0434"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "3.389.01"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.321"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "0.0₃0687"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
